# worknotes.xlsx update: "mod server for station"
# Adds new planning notes to the 全局 (sheet3 / "Global") and 路由 (sheet4 / "Routing")
# worksheets, describing the server-side station/generator module layout.

$wb  = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item(3)   # 全局
$ws4 = $wb.Worksheets.Item(4)   # 路由

# ---------------------------------------------------------------------------
# 全局 (sheet3): new dated note block starting at row 19
# ---------------------------------------------------------------------------
$ws3.Range('A19').Value = 43583
$ws3.Range('A19').NumberFormat = 'm/d/yy'
$ws3.Range('C19').Value = '数据库以电厂为单位创建'
$ws3.Range('C19').Font.Color = 255

$ws3.Range('C20').Value = '电厂内的单元创建数据表'

# ---------------------------------------------------------------------------
# 路由 (sheet4): new dated note block starting at row 17
# ---------------------------------------------------------------------------
$ws4.Range('A17').Value = 43583
$ws4.Range('A17').NumberFormat = 'm/d/yy'
$ws4.Range('C17').Value = '后端路由，作为API接口'

$ws4.Range('C21').Value = '前端路由，控制页面'
$ws4.Range('D19').Value = '所有站点的Controller相同，构造函数中，根据站点ID，选择连接对应的站点DB'
$ws4.Range('D22').Value = '页面内容划分，建议采用嵌套路由'

# ---------------------------------------------------------------------------
# back to 全局 (sheet3): object-hierarchy notes
# ---------------------------------------------------------------------------
$ws3.Range('C22').Value = '对象划分层次'
$ws3.Range('E25').Value = '电量'

$ws4.Range('D18').Value = 'controller按模块划分，如用户，设备，命名添加复数s结尾，如Users，Generators'

$ws3.Range('C27').Value = 'Server目录结构'
$ws3.Range('C28').Value = 'application\controllers'
$ws3.Range('D29').Value = 'Generators.php'
$ws3.Range('E30').Value = 'api-query_start_last_log'
$ws3.Range('E31').Value = 'api-post_start_stop_log'

$ws3.Range('C33').Value = 'application\third_party\generator\models'
$ws3.Range('D34').Value = 'Generator_model.php'
$ws3.Range('E35').Value = 'get_start_last_log'
$ws3.Range('E36').Value = 'set_start_last_log'

$ws3.Range('C38').Value = 'application\third_party\station\config'
$ws3.Range('D39').Value = 'station_config.php'

$ws3.Range('C42').Value = 'application\third_party\station\language\zh_cn'
$ws3.Range('D43').Value = 'station_lang.php'

$ws3.Range('D23').Value = '站点-在Db中体现'
$ws3.Range('E24').Value = '发电机-在db的table中，controller中体现'

# ---------------------------------------------------------------------------
# Final selections / active sheet, matching the saved workbook view state.
# ---------------------------------------------------------------------------
$ws8 = $wb.Worksheets.Item(8)   # 机组启停记录
$ws8.Range('C1:C2').Select()
$ws8.Range('C1').Activate()

$ws3.Range('G30').Select()

$ws4.Range('E27').Select()
